$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"

$ws2.Range("A1").Value = "Wert"
$ws2.Range("B1").Value = "Quadratzahl"

$row = 2
for ($n = 1; $n -le 48; $n++) {
    $ws2.Cells.Item($row, 1).Value = $n
    $ws2.Cells.Item($row, 2).Formula = "=A" + $row + "*A" + $row
    $row = $row + 1
    if ($n -eq 10) {
        $row = $row + 4
    }
}

$ws2.Range("F19").Select() | Out-Null

Write-Host "done"
